$d = $word.ActiveDocument

# List of (old, new) text replacements, in document order.
# Using MatchWholeWord=$false, MatchWildcards=$false, Forward=$true, Wrap=0 (wdFindStop)
# so each Find only matches once and we don't risk re-matching freshly
# inserted text that happens to equal an earlier "old" value.
$pairs = @(
    @("84÷3=28, 0", "83÷8=10, 3"),
    @("23÷4=5, 3", "19÷7=2, 5"),
    @("11÷4=2, 3", "92÷3=30, 2"),
    @("52÷2=26, 0", "83÷6=13, 5"),
    @("43÷6=7, 1", "17÷2=8, 1"),
    @("22÷6=3, 4", "72÷6=12, 0"),
    @("19÷2=9, 1", "93÷3=31, 0"),
    @("78÷2=39, 0", "80÷4=20, 0"),
    @("47÷3=15, 2", "81÷7=11, 4"),
    @("63÷9=7, 0", "76÷5=15, 1"),
    @("61÷9=6, 7", "47÷4=11, 3"),
    @("90÷4=22, 2", "11÷4=2, 3"),
    @("41÷3=13, 2", "88÷8=11, 0"),
    @("23÷2=11, 1", "84÷5=16, 4"),
    @("54÷6=9, 0", "16÷4=4, 0"),
    @("95÷9=10, 5", "70÷9=7, 7"),
    @("97÷9=10, 7", "48÷6=8, 0"),
    @("76÷7=10, 6", "57÷9=6, 3"),
    @("34÷7=4, 6", "17÷3=5, 2"),
    @("74÷6=12, 2", "54÷2=27, 0"),
    @("47÷8=5, 7", "96÷8=12, 0"),
    @("14÷3=4, 2", "38÷4=9, 2"),
    @("75÷8=9, 3", "55÷9=6, 1"),
    @("99÷8=12, 3", "54÷9=6, 0"),
    @("74÷4=18, 2", "89÷3=29, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2)
}
